$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 586 ("「命じられたように正しき道を辿れ」...") which shifts
# all subsequent rows up by one.
$ws.Rows.Item(586).Delete()
